$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values such as
# "312.30" or "43.467.55" are preserved exactly (Excel would otherwise
# auto-convert them to numbers and strip formatting / trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '43.467.55'
$ws.Range("E2").Value = '  +2.75%  '

# Row 3
$ws.Range("D3").Value = '2.357.25'
$ws.Range("E3").Value = '  +6.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").Value = '312.30'
$ws.Range("E5").Value = '  +5.36%  '

# Row 6
$ws.Range("D6").Value = '109.68'
$ws.Range("E6").Value = '  +0.26%  '

# Row 7
$ws.Range("D7").Value = '0.643'
$ws.Range("E7").Value = '  +2.99%  '

# Row 8
$ws.Range("E8").Value = '  -0.28%  '

# Row 9
$ws.Range("D9").Value = '0.641'
$ws.Range("E9").Value = '  +6.38%  '

# Row 10
$ws.Range("D10").Value = '43.24'
$ws.Range("E10").Value = '  -1.64%  '

# Row 11
$ws.Range("E11").Value = '  +2.91%  '

# Row 12
$ws.Range("D12").Value = '8.85'
$ws.Range("E12").Value = '  +0.68%  '

# Row 13
$ws.Range("E13").Value = '  +4.08%  '

# Row 14
$ws.Range("E14").Value = '  +2.43%  '

# Row 15
$ws.Range("D15").Value = '16.43'
$ws.Range("E15").Value = '  +8.99%  '

# Row 16
$ws.Range("D16").Value = '2.711.89'
$ws.Range("E16").Value = '  +6.19%  '

# Row 17
$ws.Range("D17").Value = '2.423.21'
$ws.Range("E17").Value = '  +9.01%  '

# Row 18
$ws.Range("D18").Value = '43.445.99'
$ws.Range("E18").Value = '  +2.56%  '

# Row 19
$ws.Range("E19").Value = '  +3.64%  '

# Row 20
$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  -1.26%  '

# Row 21
$ws.Range("D21").Value = '74.96'
$ws.Range("E21").Value = '  +3.65%  '

# Row 22
$ws.Range("E22").Value = '  -0.38%  '

# Row 23
$ws.Range("D23").Value = '2.55'
$ws.Range("E23").Value = '  +10.36%  '

# Row 24
$ws.Range("D24").Value = '256.71'
$ws.Range("E24").Value = '  +12.20%  '

# Row 25
$ws.Range("E25").Value = '  +0.69%  '

# Row 26
$ws.Range("D26").Value = '12.12'
$ws.Range("E26").Value = '  +3.70%  '

# Row 27
$ws.Range("E27").Value = '  +0.08%  '

# Row 28
$ws.Range("D28").Value = '39.27'
$ws.Range("E28").Value = '  +2.80%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '22.58'
$ws.Range("E29").Value = '  +7.57%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.17'
$ws.Range("E30").Value = '  -2.72%  '

# Row 31
$ws.Range("E31").Value = '  -0.37%  '

# Row 32
$ws.Range("D32").Value = '173.44'
$ws.Range("E32").Value = '  -0.31%  '

# Row 33
$ws.Range("D33").Value = '0.0931'
$ws.Range("E33").Value = '  +4.12%  '

# Row 34
$ws.Range("D34").Value = '6.05'
$ws.Range("E34").Value = '  +6.37%  '

# Row 35
$ws.Range("E35").Value = '  +6.00%  '

# Row 36
$ws.Range("D36").Value = '4.99'
$ws.Range("E36").Value = '  -2.02%  '

# Row 37
$ws.Range("D37").Value = '4.16'
$ws.Range("E37").Value = '  -4.57%  '

# Row 38
$ws.Range("E38").Value = '  -1.28%  '

# Row 39
$ws.Range("E39").Value = '  -0.28%  '

# Row 40
$ws.Range("D40").Value = '2.79'
$ws.Range("E40").Value = '  +15.11%  '

# Row 41
$ws.Range("D41").Value = '72.64'
$ws.Range("E41").Value = '  +0.97%  '

# Row 42
$ws.Range("D42").Value = '1.50'
$ws.Range("E42").Value = '  +14.47%  '

# Row 43
$ws.Range("E43").Value = '  -0.02%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.02%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '12.76'
$ws.Range("E45").Value = '  +1.30%  '

# Row 46
$ws.Range("E46").Value = '  +3.89%  '

# Row 47
$ws.Range("D47").Value = '9.37'
$ws.Range("E47").Value = '  +11.09%  '

# Row 48
$ws.Range("D48").Value = '111.37'
$ws.Range("E48").Value = '  +7.63%  '

# Row 49
$ws.Range("E49").Value = '  +0.35%  '

# Row 50
$ws.Range("E50").Value = '  +2.84%  '

# Row 51
$ws.Range("D51").Value = '0.472'
$ws.Range("E51").Value = '  +7.65%  '
